$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Sprint 7 story statuses ---
# Row 23 ("I want include options to mark who paid for the expense") moves
# from NOT STARTED -> IN PROGRESS. Grab C22's current "Neutral" (with
# border) formatting first -- that's exactly the look every other
# "IN PROGRESS" cell uses -- before C22 itself gets overwritten below.
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null

# Row 22 ("I want to include ads") moves from IN PROGRESS -> DONE. Reuse
# the "Good" (with border) formatting already used by every other DONE
# cell (e.g. C2).
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Apply the new status text now that formats are sorted out.
$ws.Range("C22").Value = "DONE"
$ws.Range("C23").Value = "IN PROGRESS"

# --- Update the active-cell selection left behind by the edit ---
$ws.Range("G19").Select() | Out-Null
